$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "60.698.18"
$ws.Range("D3").Value = "2.906.52"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.41%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "2.915.71"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("E10").Value = "  -4.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.359"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").Value = "3.416.23"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "60.692.28"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "2.910.88"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.181"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.81%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.09%  "
$ws.Range("D30").Value = "0.0₃0849"
$ws.Range("E30").Value = "  -9.71%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.02%  "
$ws.Range("E36").Value = "  -6.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.30%  "
$ws.Range("E38").Value = "  -4.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  -4.51%  "
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("D42").Value = "2.295.86"
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0586"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0925"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "249.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.46%  "
